$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing "REQ 6`" header in H1 (value cleared, style kept)
$ws.Range("H1").Value = $null

# Fill in the new Totals row values (row 2)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

# Add Totals column values for TC4 and TC5 rows
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1

# Mark the new intersections with "X"
$ws.Range("E5").Value = "X"
$ws.Range("F6").Value = "X"
$ws.Range("G7").Value = "X"

# Remove row 8 (TC6) entirely - clear label and any formatting content
$ws.Range("A8").Value = $null
